$wb = $excel.ActiveWorkbook

# --- Sheet "Logs": append the new log row (row 24) ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A24").Value = "Opvolging afspraak"
$logs.Range("B24").Value = "mailmind.test@zohomail.eu"
$logs.Range("D24").Value = "Planning / Afspraak"
$logs.Range("F24").Value = "2025-08-28 21:12:24"
$logs.Range("G24").Value = "Nee"
$logs.Range("H24").Value = "Ja"
$logs.Range("I24").Value = "Nee"
$logs.Range("J24").Value = "Nee"

# Extend the conditional formatting ranges from row 23 to row 24 so they
# keep covering the whole data range (D/G/H/I/J columns).
$dCond = $logs.Range("D2:D23").FormatConditions
for ($i = 1; $i -le $dCond.Count; $i++) {
    $dCond.Item($i).ModifyAppliesToRange($logs.Range("D2:D24"))
}

$gCond = $logs.Range("G2:G23").FormatConditions
for ($i = 1; $i -le $gCond.Count; $i++) {
    $gCond.Item($i).ModifyAppliesToRange($logs.Range("G2:G24"))
}

$hCond = $logs.Range("H2:H23").FormatConditions
for ($i = 1; $i -le $hCond.Count; $i++) {
    $hCond.Item($i).ModifyAppliesToRange($logs.Range("H2:H24"))
}

$iCond = $logs.Range("I2:I23").FormatConditions
for ($i = 1; $i -le $iCond.Count; $i++) {
    $iCond.Item($i).ModifyAppliesToRange($logs.Range("I2:I24"))
}

$jCond = $logs.Range("J2:J23").FormatConditions
for ($i = 1; $i -le $jCond.Count; $i++) {
    $jCond.Item($i).ModifyAppliesToRange($logs.Range("J2:J24"))
}

# --- Sheet "Dashboard": refresh the category summary table ---
# "Planning / Afspraak" now has 2 occurrences and moves above "Overig" (1).
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A4").Value = "Planning / Afspraak"
$dash.Range("B4").Value = 2
$dash.Range("A5").Value = "Overig"
$dash.Range("B5").Value = 1
